# Auto-generated edit script: update Maduin_Profits market-price derived cells
# Mirrors a scheduled market-data refresh across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 59
$ws.Range("H59").Value = 4500
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 4500
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 13500
$ws.Range("N59").Value = -14614
$ws.Range("M59").ClearContents()

# Row 64
$ws.Range("H64").Value = 9544.727999999999
$ws.Range("I64").Value = 5624.75
$ws.Range("J64").Value = 19998
$ws.Range("K64").Value = 5624.75
$ws.Range("L64").Value = 19998
$ws.Range("M64").Value = -5376.75
$ws.Range("N64").Value = -20494

# Row 67
$ws.Range("H67").Value = 9544.727999999999
$ws.Range("I67").Value = 5624.75
$ws.Range("J67").Value = 19998
$ws.Range("K67").Value = 5624.75
$ws.Range("L67").Value = 19998
$ws.Range("M67").Value = -4766.75
$ws.Range("N67").Value = -21714

# Row 86
$ws.Range("H86").Value = 8399.223
$ws.Range("I86").Value = 7720
$ws.Range("J86").Value = 9248.25
$ws.Range("K86").Value = 7720
$ws.Range("L86").Value = 9248.25
$ws.Range("M86").Value = -6597
$ws.Range("N86").Value = -11494.25

# Row 89
$ws.Range("H89").Value = 8399.223
$ws.Range("I89").Value = 7720
$ws.Range("J89").Value = 9248.25
$ws.Range("K89").Value = 38600
$ws.Range("L89").Value = 46241.25
$ws.Range("M89").Value = -32984
$ws.Range("N89").Value = -57473.25

# Row 132
$ws.Range("H132").Value = 2392.4482
$ws.Range("I132").Value = 2341.25
$ws.Range("J132").Value = 2506.2222
$ws.Range("K132").Value = 7023.75
$ws.Range("L132").Value = 7518.6666
$ws.Range("M132").Value = -4493.75
$ws.Range("N132").Value = -12578.6666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 56942.332
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 56942.332
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 56942.332
$ws.Range("N24").Value = -57690.332

# Row 45
$ws.Range("H45").Value = 3551.3
$ws.Range("I45").Value = 1171
$ws.Range("J45").Value = 4571.4287
$ws.Range("K45").Value = 1171
$ws.Range("L45").Value = 4571.4287
$ws.Range("M45").Value = -794
$ws.Range("N45").Value = -5325.4287

# Row 61
$ws.Range("H61").Value = 1087.2
$ws.Range("I61").Value = 984
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 984
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -772
$ws.Range("N61").Value = -1924

# Row 97
$ws.Range("H97").Value = 1208.1666
$ws.Range("I97").Value = 695.8461
$ws.Range("J97").Value = 2540.2
$ws.Range("K97").Value = 695.8461
$ws.Range("L97").Value = 2540.2
$ws.Range("M97").Value = -199.8461
$ws.Range("N97").Value = -3532.2

# Row 100
$ws.Range("H100").Value = 56942.332
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 56942.332
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 56942.332
$ws.Range("N100").Value = -59106.332

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()

# Row 110
$ws.Range("H110").Value = 970
$ws.Range("I110").Value = 970
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 970
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1075

# Row 136
$ws.Range("H136").Value = 1087.2
$ws.Range("I136").Value = 984
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 2952
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -402
$ws.Range("N136").Value = -9600

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 3710.4614
$ws.Range("I99").Value = 2913.6
$ws.Range("J99").Value = 6366.6665
$ws.Range("K99").Value = 2913.6
$ws.Range("L99").Value = 6366.6665
$ws.Range("M99").Value = -1415.6
$ws.Range("N99").Value = -9362.666499999999

# Row 105
$ws.Range("H105").Value = 3668.625
$ws.Range("I105").Value = 3323
$ws.Range("J105").Value = 5166.3335
$ws.Range("K105").Value = 3323
$ws.Range("L105").Value = 5166.3335
$ws.Range("M105").Value = -1576
$ws.Range("N105").Value = -8660.333500000001

# Row 107
$ws.Range("H107").Value = 1199
$ws.Range("I107").Value = 1199
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1199
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 721
$ws.Range("N107").ClearContents()

# Row 134
$ws.Range("H134").Value = 1563.4615
$ws.Range("I134").Value = 1443.75
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 4331.25
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -1796.25
$ws.Range("N134").Value = -14070

# Row 137
$ws.Range("H137").Value = 64999
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 64999
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 64999
$ws.Range("N137").Value = -75199

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 972.1667
$ws.Range("I16").Value = 710
$ws.Range("J16").Value = 1103.25
$ws.Range("K16").Value = 710
$ws.Range("L16").Value = 1103.25
$ws.Range("M16").Value = -423
$ws.Range("N16").Value = -1677.25

# Row 113
$ws.Range("H113").Value = 972.1667
$ws.Range("I113").Value = 710
$ws.Range("J113").Value = 1103.25
$ws.Range("K113").Value = 710
$ws.Range("L113").Value = 1103.25
$ws.Range("M113").Value = 1460
$ws.Range("N113").Value = -5443.25

# Row 122
$ws.Range("H122").Value = 1338.3334
$ws.Range("I122").Value = 600
$ws.Range("J122").Value = 1707.5
$ws.Range("K122").Value = 1800
$ws.Range("L122").Value = 5122.5
$ws.Range("M122").Value = 650
$ws.Range("N122").Value = -10022.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 10713.857
$ws.Range("I80").Value = 2499.5
$ws.Range("J80").Value = 13999.6
$ws.Range("K80").Value = 7498.5
$ws.Range("L80").Value = 41998.8
$ws.Range("M80").Value = -6562.5
$ws.Range("N80").Value = -43870.8

# Row 83
$ws.Range("H83").Value = 10713.857
$ws.Range("I83").Value = 2499.5
$ws.Range("J83").Value = 13999.6
$ws.Range("K83").Value = 22495.5
$ws.Range("L83").Value = 125996.4
$ws.Range("M83").Value = -17815.5
$ws.Range("N83").Value = -135356.4

# Row 107
$ws.Range("H107").Value = 1068.5555
$ws.Range("I107").Value = 862.5
$ws.Range("J107").Value = 1233.4
$ws.Range("K107").Value = 2587.5
$ws.Range("L107").Value = 3700.2
$ws.Range("M107").Value = -667.5
$ws.Range("N107").Value = -7540.200000000001

# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# Row 128
$ws.Range("H128").Value = 278506
$ws.Range("I128").Value = 278506
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 835518
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -830538

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1495.8334
$ws.Range("I102").Value = 1195.2
$ws.Range("J102").Value = 2999
$ws.Range("K102").Value = 1195.2
$ws.Range("L102").Value = 2999
$ws.Range("M102").Value = 426.8
$ws.Range("N102").Value = -6243

# Row 126
$ws.Range("H126").Value = 4324
$ws.Range("I126").Value = 4049.5
$ws.Range("J126").Value = 4598.5
$ws.Range("K126").Value = 12148.5
$ws.Range("L126").Value = 13795.5
$ws.Range("M126").Value = -9678.5
$ws.Range("N126").Value = -18735.5

# Row 132
$ws.Range("H132").Value = 7560.857
$ws.Range("I132").Value = 7587.4
$ws.Range("J132").Value = 7494.5
$ws.Range("K132").Value = 22762.2
$ws.Range("L132").Value = 22483.5
$ws.Range("M132").Value = -20232.2
$ws.Range("N132").Value = -27543.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 190.1
$ws.Range("I16").Value = 190.1
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 190.1
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -20.09999999999999
$ws.Range("N16").ClearContents()

# Row 22
$ws.Range("H22").Value = 924.3333
$ws.Range("I22").Value = 924.3333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 924.3333
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -629.3333
$ws.Range("N22").ClearContents()

# Row 27
$ws.Range("H27").Value = 924.3333
$ws.Range("I27").Value = 924.3333
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 924.3333
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -817.3333
$ws.Range("N27").ClearContents()

# Row 40
$ws.Range("H40").Value = 5320.25
$ws.Range("I40").Value = 5417.375
$ws.Range("J40").Value = 5126
$ws.Range("K40").Value = 5417.375
$ws.Range("L40").Value = 5126
$ws.Range("M40").Value = -5281.375
$ws.Range("N40").Value = -5398

# Row 43
$ws.Range("H43").Value = 18507
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 18507
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 18507
$ws.Range("N43").Value = -18893

# Row 46
$ws.Range("H46").Value = 1663.76
$ws.Range("I46").Value = 1185
$ws.Range("J46").Value = 2273.0908
$ws.Range("K46").Value = 1185
$ws.Range("L46").Value = 2273.0908
$ws.Range("M46").Value = -997
$ws.Range("N46").Value = -2649.0908

# Row 93
$ws.Range("H93").Value = 2250
$ws.Range("I93").Value = 2000
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 2000
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -752
$ws.Range("N93").Value = -4996
